# Adds a "Player Info" sheet in front of the existing sheets and switches
# the MATCH_CARD_LINK columns (full howstat.com URLs) on the ODI Batting /
# ODI Bowling sheets to a plain MATCH_CODE (just the numeric match id).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Player Info" sheet as the very first sheet.
# ---------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").Value = "'4656"
$playerInfo.Range("B2").Value = "Navdeep Amarjeet Saini"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast"

# ---------------------------------------------------------------------
# 2. "ODI Batting" sheet: MATCH_CARD_LINK (col D) -> MATCH_CODE
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @("4388", "4399", "4400", "4406", "4410", "4435", "4436", "4485")
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    $batting.Cells.Item($row, 4).Value = "'" + $battingCodes[$i]
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling" sheet: MATCH_CARD_LINK (col B) -> MATCH_CODE
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @("4388", "4399", "4400", "4406", "4410", "4435", "4436", "4485")
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $row = $i + 2
    $bowling.Cells.Item($row, 2).Value = "'" + $bowlingCodes[$i]
}
